# Insert a new data row at row 296 (pushing the existing rows 296-378 down
# to 297-379) and populate it with the new record described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 296; this shifts all rows from
# 296 downward to 297 onward, preserving their data/formatting.
$ws.Rows.Item(296).Insert()

# Populate the newly inserted (now empty) row 296 with the new record.
$ws.Cells.Item(296, 1).Value2  = 9
$ws.Cells.Item(296, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(296, 3).Value2  = "Metropolitana"
$ws.Cells.Item(296, 4).Value2  = 44642
$ws.Cells.Item(296, 5).Value2  = 13
$ws.Cells.Item(296, 6).Value2  = "Fruta"
$ws.Cells.Item(296, 7).Value2  = 100108
$ws.Cells.Item(296, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(296, 9).Value2  = 100108002
$ws.Cells.Item(296, 10).Value2 = "Mango"
$ws.Cells.Item(296, 11).Value2 = "Sin especificar"
$ws.Cells.Item(296, 12).Value2 = "Primera"
$ws.Cells.Item(296, 13).Value2 = 680
$ws.Cells.Item(296, 14).Value2 = 6500
$ws.Cells.Item(296, 15).Value2 = 7000
$ws.Cells.Item(296, 16).Value2 = 6743
$ws.Cells.Item(296, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(296, 18).Value2 = "Perú"
$ws.Cells.Item(296, 19).Value2 = 1686
$ws.Cells.Item(296, 20).Value2 = 4
